$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "52.154.11"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.990.94"
$ws.Range("E3").Value = "  +0.53%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "354.31"
$ws.Range("E5").Value = "  +0.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "107.89"
$ws.Range("E6").Value = "  -4.29%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.562"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.625"
$ws.Range("E9").Value = "  -1.15%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.38"
$ws.Range("E10").Value = "  -3.47%  "
$ws.Range("E11").Value = "  +2.37%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0858"
$ws.Range("E12").Value = "  -4.73%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.34"
$ws.Range("E13").Value = "  -3.33%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.463.18"
$ws.Range("E14").Value = "  +0.39%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.76"
$ws.Range("E15").Value = "  -2.05%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.989.60"
$ws.Range("E16").Value = "  +0.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.02"
$ws.Range("E17").Value = "  +3.28%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "52.229.78"
$ws.Range("E18").Value = "  +0.14%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.51"
$ws.Range("E19").Value = "  +5.78%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.56"
$ws.Range("E20").Value = "  -1.82%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.67"
$ws.Range("E21").Value = "  -5.20%  "
$ws.Range("E22").Value = "  -1.69%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.51"
$ws.Range("E23").Value = "  -2.64%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "264.13"
$ws.Range("E24").Value = "  -2.38%  "
$ws.Range("E25").Value = "  -2.30%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.182"
$ws.Range("E26").Value = "  +0.24%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.89"
$ws.Range("E27").Value = "  -2.19%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.57"
$ws.Range("E28").Value = "  -1.99%  "
$ws.Range("E29").Value = "  -0.26%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.108"
$ws.Range("E30").Value = "  -5.66%  "
$ws.Range("E31").Value = "  -3.83%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.37"
$ws.Range("E32").Value = "  +1.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "36.70"
$ws.Range("E33").Value = "  -2.30%  "
$ws.Range("E34").Value = "  +7.43%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "50.90"
$ws.Range("E35").Value = "  -4.04%  "
$ws.Range("E36").Value = "  -1.18%  "
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.20"
$ws.Range("E38").Value = "  -6.14%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "17.95"
$ws.Range("E39").Value = "  -5.65%  "
$ws.Range("E40").Value = "  -4.95%  "
$ws.Range("E41").Value = "  +0.49%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.118"
$ws.Range("E42").Value = "  -0.19%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.93"
$ws.Range("E43").Value = "  -3.84%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "122.33"
$ws.Range("E44").Value = "  +6.68%  "
$ws.Range("E45").Value = "  -1.39%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.125.54"
$ws.Range("E46").Value = "  -2.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.39"
$ws.Range("E47").Value = "  -4.60%  "
$ws.Range("E48").Value = "  -5.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.250"
$ws.Range("E49").Value = "  +1.93%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0335"
$ws.Range("E50").Value = "  -2.53%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.923"
$ws.Range("E51").Value = "  -2.26%  "
